# Trade #2 closed at 2026-02-16 21:50:46 - leadlag UP +0.000%
#
# Append a new "OPEN" trade row (row 3) to both the "All Trades" and the
# "leadlag" worksheets, mirroring the layout already used for trade #1 in
# row 2 of each sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "leadlag")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(3, 1).Value = 2
    # Leading apostrophe forces the date-shaped string to be stored as text
    # (matching column B's existing "Date" text cells) instead of being
    # auto-converted to a date serial number.
    $ws.Cells.Item(3, 2).Value = "'2026-02-16"
    $ws.Cells.Item(3, 3).Value = "21:50:46"
    $ws.Cells.Item(3, 4).Value = "leadlag"
    $ws.Cells.Item(3, 5).Value = "UP"
    $ws.Cells.Item(3, 6).Value = 68350.87
    $ws.Cells.Item(3, 8).Value = "OPEN"
    $ws.Cells.Item(3, 9).Value = 0
    $ws.Cells.Item(3, 10).Value = 0
    $ws.Cells.Item(3, 11).Value = 100
    $ws.Cells.Item(3, 12).Value = 0.75
    $ws.Cells.Item(3, 13).Value = "Coinbase leading with 0.097% move"
    $ws.Cells.Item(3, 15).Value = 0
}
